# Insert a new row 455 (shifts existing rows 455-498 down to 456-499)
# and populate it with the new observation, matching the committed diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(455).Insert()

$row = 455
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 45106
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100112017
$ws.Cells.Item($row, 7).Value = "Apio"
$ws.Cells.Item($row, 8).Value = "Americana (o)"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 210
$ws.Cells.Item($row, 11).Value = 8000
$ws.Cells.Item($row, 12).Value = 8000
$ws.Cells.Item($row, 13).Value = 8000
$ws.Cells.Item($row, 14).Value = "$/docena de matas"
$ws.Cells.Item($row, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 16).Value = 1333
$ws.Cells.Item($row, 17).Value = 6
$ws.Cells.Item($row, 18).Value = "Hortaliza"
